$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 -> becomes what was row 7
$ws.Range("D4").Value = 44719
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14400
$ws.Range("Q4").Value = "$/caja 18 kilos granel"
$ws.Range("R4").Value = "Región del Maule"
$ws.Range("S4").Value = 800
$ws.Range("T4").Value = 18

# Row 6 -> becomes what was row 8
$ws.Range("D6").Value = 44714
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 806
$ws.Range("T6").Value = 18

# Row 7 -> becomes what was row 6
$ws.Range("D7").Value = 44334
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11500
$ws.Range("Q7").Value = "$/caja 12 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 11500
$ws.Range("T7").Value = 1

# Row 8 -> becomes what was row 4
$ws.Range("D8").Value = 44708
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 13000
$ws.Range("P8").Value = 12571
$ws.Range("Q8").Value = "$/caja 12 kilos empedrada"
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 1048
$ws.Range("T8").Value = 12
